$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.199.34"
$ws.Range("E2").Value = "  +2.31%  "
$ws.Range("D3").Value = "3.698.39"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "582.62"
$ws.Range("E5").Value = "  -0.30%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "178.33"
$ws.Range("E6").Value = "  +0.78%  "
$ws.Range("D7").Value = "3.687.28"
$ws.Range("E7").Value = "  +7.76%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.618"
$ws.Range("E8").Value = "  +4.10%  "
$ws.Range("E9").Value = "  -0.02%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.199"
$ws.Range("E10").Value = "  -0.63%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.97"
$ws.Range("E11").Value = "  +28.25%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.612"
$ws.Range("E12").Value = "  +4.82%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "49.21"
$ws.Range("E13").Value = "  +0.87%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000288"
$ws.Range("E14").Value = "  +2.06%  "
$ws.Range("D15").Value = "4.291.96"
$ws.Range("E15").Value = "  +7.92%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "680.30"
$ws.Range("E16").Value = "  -1.72%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "9.03"
$ws.Range("E17").Value = "  +4.66%  "
$ws.Range("D18").Value = "3.704.35"
$ws.Range("E18").Value = "  +8.07%  "
$ws.Range("D19").Value = "71.314.72"
$ws.Range("E19").Value = "  +2.40%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.123"
$ws.Range("E20").Value = "  +0.95%  "
$ws.Range("E21").Value = "  +1.74%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.62"
$ws.Range("E22").Value = "  +2.28%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.945"
$ws.Range("E23").Value = "  +5.38%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "17.44"
$ws.Range("E24").Value = "  +2.80%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "102.31"
$ws.Range("E25").Value = "  +1.02%  "
$ws.Range("E26").Value = "  +1.88%  "
$ws.Range("E27").Value = "  +7.00%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.28"
$ws.Range("E28").Value = "  +7.20%  "
$ws.Range("E29").Value = "  -0.04%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.15"
$ws.Range("E30").Value = "  +4.82%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.43"
$ws.Range("E31").Value = "  +5.26%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "9.17"
$ws.Range("E32").Value = "  +4.46%  "
$ws.Range("E33").Value = "  -1.49%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.55"
$ws.Range("E34").Value = "  +6.46%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.08"
$ws.Range("E35").Value = "  +10.15%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "581.22"
$ws.Range("E36").Value = "  +1.55%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.108"
$ws.Range("E38").Value = "  +4.69%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "58.69"
$ws.Range("E39").Value = "  +0.80%  "
$ws.Range("D40").Value = "3.681.40"
$ws.Range("E40").Value = "  +3.00%  "
$ws.Range("E41").Value = "  -0.16%  "
$ws.Range("E42").Value = "  +3.79%  "
$ws.Range("B43").Value = "TheGraph"
$ws.Range("C43").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.353"
$ws.Range("E43").Value = "  +6.31%  "
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0458"
$ws.Range("E44").Value = "  +9.80%  "
$ws.Range("D45").Value = "0.0₃0769"
$ws.Range("E45").Value = "  +4.97%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "35.73"
$ws.Range("E46").Value = "  +2.23%  "
$ws.Range("E47").Value = "  +4.06%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.94"
$ws.Range("E48").Value = "  +11.02%  "
$ws.Range("E49").Value = "  +4.21%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "135.14"
$ws.Range("E50").Value = "  +1.86%  "
$ws.Range("E51").Value = "  +11.31%  "
